{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst replacements = [\n  [\n    \"Play Lucky Halloween Free: Review of Red Tiger Game\",\n    \"Play Lucky Halloween Slot and Win Big for Free!\",\n  ],\n  [\"Entertaining Halloween theme\", \"Spooky Halloween theme\"],\n  [\"20 paylines for winning combinations\", \"Variety of special features\"],\n  [\"Special features like wilds and multipliers\", \"Engaging gameplay\"],\n  [\n    \"Eerie sound effects adding to the spooky atmosphere\",\n    \"Opportunities to win big\",\n  ],\n  [\n    \"May not be suitable for those who are easily frightened\",\n    \"May not be suitable for those easily frightened\",\n  ],\n  [\"No progressive jackpots\", \"Limited betting range\"],\n  [\n    \"Spooky Lucky Halloween slot game from Red Tiger. Review includes gameplay, special features, symbols, and sound effects. Play Lucky Halloween free or for real money.\",\n    \"Play Lucky Halloween slot game for free and enjoy the spooky Halloween theme and chances to win big.\",\n  ],\n];\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d (ActiveDocument) are pre-seeded.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Lucky Halloween Free: Review of Red Tiger Game\", \"Play Lucky Halloween Slot and Win Big for Free!\"),\n    @(\"Entertaining Halloween theme\", \"Spooky Halloween theme\"),\n    @(\"20 paylines for winning combinations\", \"Variety of special features\"),\n    @(\"Special features like wilds and multipliers\", \"Engaging gameplay\"),\n    @(\"Eerie sound effects adding to the spooky atmosphere\", \"Opportunities to win big\"),\n    @(\"May not be suitable for those who are easily frightened\", \"May not be suitable for those easily frightened\"),\n    @(\"No progressive jackpots\", \"Limited betting range\"),\n    @(\"Spooky Lucky Halloween slot game from Red Tiger. Review includes gameplay, special features, symbols, and sound effects. Play Lucky Halloween free or for real money.\", \"Play Lucky Halloween slot game for free and enjoy the spooky Halloween theme and chances to win big.\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)  # wdReplaceAll\n}\n"}
